$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (template header row)
$ws.Range("B1").Value = "Title1"
$ws.Range("C1").Value = "+Tag1"
$ws.Range("D1").Value = "+Tag2"

# Row 2
$ws.Range("B2").Value = "Title1"
$ws.Range("C2").Value = "+Tag1"
$ws.Range("D2").Value = "+Tag2"
$ws.Range("E2").ClearContents()

# Row 3
$ws.Range("B3").Value = "Title2"
$ws.Range("C3").Value = "+Tag1"
$ws.Range("D3").Value = "+Tag2"
$ws.Range("E3").Value = "?"

# Row 4 (advanced main analysis: allows parameter changes)
$ws.Range("B4").Value = "Title3"
$ws.Range("C4").Value = "+ALL"
$ws.Range("D4").Value = "`$mask.nii"

# Row 5
$ws.Range("A5").Value = "B"
$ws.Range("B5").Value = "Title4"
$ws.Range("C5").Value = "+Tag1"
$ws.Range("D5").Value = "-Tag2"
$ws.Range("E5").ClearContents()

# Row 6
$ws.Range("B6").Value = "Title1"
$ws.Range("C6").Value = "+Tag1"
$ws.Range("D6").Value = "+Tag2"

# Update active selection to reflect where the author left off editing
$ws.Range("A5").Select()
